$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws1.Range("H15").Value = 2795.8262
$ws1.Range("I15").Value = 2795.8262
$ws1.Range("K15").Value = 8387.4786
$ws1.Range("M15").Value = -8218.4786
$ws1.Range("H40").Value = 2111.111
$ws1.Range("I40").Value = 1666.6666
$ws1.Range("J40").Value = 2333.3333
$ws1.Range("K40").Value = 1666.6666
$ws1.Range("L40").Value = 2333.3333
$ws1.Range("M40").Value = -1491.6666
$ws1.Range("N40").Value = -2683.3333
$ws1.Range("H43").Value = 2149.8333
$ws1.Range("I43").Value = 1500
$ws1.Range("J43").Value = 2279.8
$ws1.Range("K43").Value = 1500
$ws1.Range("L43").Value = 2279.8
$ws1.Range("M43").Value = -1431
$ws1.Range("N43").Value = -2417.8
$ws1.Range("H86").Value = 7478.1113
$ws1.Range("I86").Value = 1582.2222
$ws1.Range("J86").Value = 13374
$ws1.Range("K86").Value = 1582.2222
$ws1.Range("L86").Value = 13374
$ws1.Range("M86").Value = -459.2221999999999
$ws1.Range("N86").Value = -15620
$ws1.Range("H89").Value = 7478.1113
$ws1.Range("I89").Value = 1582.2222
$ws1.Range("J89").Value = 13374
$ws1.Range("K89").Value = 7911.111
$ws1.Range("L89").Value = 66870
$ws1.Range("M89").Value = -2295.111
$ws1.Range("N89").Value = -78102
$ws1.Range("H125").Value = 1500
$ws1.Range("I125").Value = 1500
$ws1.Range("K125").Value = 13500
$ws1.Range("M125").Value = -11040
$ws1.Range("H129").Value = 1211.1831
$ws1.Range("J129").Value = 1221.3429
$ws1.Range("L129").Value = 3664.0287
$ws1.Range("N129").Value = -13664.0287
$ws1.Range("H135").Value = 7523.75
$ws1.Range("I135").Value = 734.6111
$ws1.Range("J135").Value = 27891.166
$ws1.Range("K135").Value = 6611.4999
$ws1.Range("L135").Value = 251020.494
$ws1.Range("M135").Value = -4076.4999
$ws1.Range("N135").Value = -256090.494
$ws2 = $wb.Worksheets.Item("ARM")
$ws2.Range("H36").Value = 29999.5
$ws2.Range("I36").Value = 29999
$ws2.Range("J36").Value = 30000
$ws2.Range("K36").Value = 29999
$ws2.Range("L36").Value = 30000
$ws2.Range("M36").Value = -29653
$ws2.Range("N36").Value = -30692
$ws2.Range("H37").Value = 15995
$ws2.Range("I37").Value = 2000
$ws2.Range("K37").Value = 2000
$ws2.Range("M37").Value = -1727
$ws2.Range("H102").Value = 1304.6471
$ws2.Range("I102").Value = 1090.7693
$ws2.Range("K102").Value = 1090.7693
$ws2.Range("M102").Value = 531.2307000000001
$ws2.Range("H132").Value = 23195.125
$ws2.Range("I132").Value = 2233.55
$ws2.Range("J132").Value = 128003
$ws2.Range("K132").Value = 6700.650000000001
$ws2.Range("L132").Value = 384009
$ws2.Range("M132").Value = -4170.650000000001
$ws2.Range("N132").Value = -389069
$ws3 = $wb.Worksheets.Item("BSM")
$ws3.Range("H99").Value = 1733.8334
$ws3.Range("I99").Value = 1487.375
$ws3.Range("J99").Value = 1931
$ws3.Range("K99").Value = 1487.375
$ws3.Range("L99").Value = 1931
$ws3.Range("M99").Value = 10.625
$ws3.Range("N99").Value = -4927
$ws4 = $wb.Worksheets.Item("CRP")
$ws4.Range("H48").Value = 0
$ws4.Range("J48").Value = 0
$ws4.Range("L48").Value = 0
$ws4.Range("N48").Value = ""
$ws4.Range("H122").Value = 1149.8422
$ws4.Range("I122").Value = 868.4286
$ws4.Range("K122").Value = 2605.2858
$ws4.Range("M122").Value = -155.2857999999997
$ws5 = $wb.Worksheets.Item("CUL")
$ws5.Range("H5").Value = 1728.5454
$ws5.Range("I5").Value = 1445.4445
$ws5.Range("K5").Value = 4336.333500000001
$ws5.Range("M5").Value = -4224.333500000001
$ws5.Range("H126").Value = 3116
$ws5.Range("J126").Value = 6500
$ws5.Range("L126").Value = 19500
$ws5.Range("N126").Value = -29380
$ws5.Range("H131").Value = 794.25
$ws5.Range("J131").Value = 798.2143
$ws5.Range("L131").Value = 2394.6429
$ws5.Range("N131").Value = -12474.6429
$ws5.Range("H135").Value = 1728.5454
$ws5.Range("I135").Value = 1445.4445
$ws5.Range("K135").Value = 13009.0005
$ws5.Range("M135").Value = -10474.0005
$ws5.Range("H137").Value = 1943.5625
$ws5.Range("I137").Value = 519.1667
$ws5.Range("J137").Value = 2798.2
$ws5.Range("K137").Value = 1557.5001
$ws5.Range("L137").Value = 8394.599999999999
$ws5.Range("M137").Value = 3542.4999
$ws5.Range("N137").Value = -18594.6
$ws6 = $wb.Worksheets.Item("GSM")
$ws6.Range("H46").Value = 23975
$ws6.Range("J46").Value = 23975
$ws6.Range("L46").Value = 23975
$ws6.Range("N46").Value = -24287
$ws6.Range("H107").Value = 4808798.5
$ws6.Range("I107").Value = 239.9
$ws6.Range("J107").Value = 12823062
$ws6.Range("K107").Value = 239.9
$ws6.Range("L107").Value = 12823062
$ws6.Range("M107").Value = 1680.1
$ws6.Range("N107").Value = -12826902
$ws6.Range("H122").Value = 2666.7334
$ws6.Range("I122").Value = 2171.4285
$ws6.Range("J122").Value = 3100.125
$ws6.Range("K122").Value = 6514.2855
$ws6.Range("L122").Value = 9300.375
$ws6.Range("M122").Value = -4064.2855
$ws6.Range("N122").Value = -14200.375
$ws6.Range("H131").Value = 45000
$ws6.Range("J131").Value = 45000
$ws6.Range("L131").Value = 45000
$ws6.Range("N131").Value = -55080
$ws7 = $wb.Worksheets.Item("LTW")
$ws7.Range("H7").Value = 3243.3635
$ws7.Range("I7").Value = 3081.7896
$ws7.Range("K7").Value = 3081.7896
$ws7.Range("M7").Value = -2969.7896
$ws7.Range("H40").Value = 3603.4119
$ws7.Range("I40").Value = 3230.8
$ws7.Range("K40").Value = 3230.8
$ws7.Range("M40").Value = -3094.8
$ws7.Range("H126").Value = 3243.3635
$ws7.Range("I126").Value = 3081.7896
$ws7.Range("K126").Value = 9245.3688
$ws7.Range("M126").Value = -6775.3688
$ws7.Range("H132").Value = 1962.88
$ws7.Range("I132").Value = 1390.3334
$ws7.Range("J132").Value = 2821.7
$ws7.Range("K132").Value = 4171.0002
$ws7.Range("L132").Value = 8465.099999999999
$ws7.Range("M132").Value = -1641.0002
$ws7.Range("N132").Value = -13525.1
$ws8 = $wb.Worksheets.Item("WVR")
$ws8.Range("H81").Value = 83334540
$ws8.Range("I81").Value = 1412.7
$ws8.Range("K81").Value = 2825.4
$ws8.Range("M81").Value = -1764.4
$ws8.Range("H84").Value = 83334540
$ws8.Range("I84").Value = 1412.7
$ws8.Range("K84").Value = 14127
$ws8.Range("M84").Value = -8823
$ws8.Range("H107").Value = 2274012.5
$ws8.Range("I107").Value = 1112.75
$ws8.Range("K107").Value = 3338.25
$ws8.Range("M107").Value = -1418.25
$ws8.Range("H132").Value = 1820.1
$ws8.Range("I132").Value = 1025.375
$ws8.Range("K132").Value = 3076.125
$ws8.Range("M132").Value = -546.125
$ws8.Range("H136").Value = 35716220
$ws8.Range("I136").Value = 52633372
$ws8.Range("J136").Value = 2233.6667
$ws8.Range("K136").Value = 157900116
$ws8.Range("L136").Value = 6701.000100000001
$ws8.Range("M136").Value = -157897566
$ws8.Range("N136").Value = -11801.0001
